$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before row 723, shifting existing rows 723:801 down to 728:806
$ws.Range("A723:T727").EntireRow.Insert()

# Common (unchanged) columns for this block of data
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria = "Nectarín"

function Set-Row($r, $fecha, $k, $l, $m, $n, $o, $p, $q, $rOrigen, $s, $t) {
    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $k
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $rOrigen
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}

Set-Row 723 44918 "Artic Star"    "Especial" 300 18000 18000 18000 "`$/caja 18 kilos granel" "Provincia de San Felipe de Aconcagua" 1000 18
Set-Row 724 44918 "Artic Star"    "Primera"  220 14400 14400 14400 "`$/caja 18 kilos granel" "Provincia de San Felipe de Aconcagua" 800  18
Set-Row 725 44918 "Artic Star"    "Segunda"  250 11000 11000 11000 "`$/caja 18 kilos granel" "Provincia de San Felipe de Aconcagua" 611  18
Set-Row 726 44918 "Early Diamond" "Especial" 250 19200 19200 19200 "`$/caja 16 kilos granel" "Región de O'Higgins" 1200 16
Set-Row 727 44918 "Early Diamond" "Primera"  300 16000 16000 16000 "`$/caja 16 kilos granel" "Región de O'Higgins" 1000 16

# Note: the D column (date) already keeps the existing date number format (style index 2)
# because Range.Insert() (shift-down) propagates the format from the row above the
# insertion point (the former row 723) into the newly inserted rows 723:727.
